$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# Title
Replace-Text "Quantum Computing: A Paradigm Shift in Computational Power" "History: A Portal to the Past"

# Author line (merges " Angela Q" + "." + " Thompson" runs into one)
Replace-Text " Angela Q. Thompson" " Rebecca Mackenzie"

# Email / affiliation line (merges "thompson@qbit" + "." + "institute" runs into one)
Replace-Text "thompson@qbit.institute" "at"

# Body paragraph, sentence by sentence
Replace-Text "In the realm of computational science, a transformative technology known as quantum computing is poised to revolutionize our understanding of information processing and problem-solving" "History is a vast and captivating subject that holds the power to transport us through time, unveil the mysteries of the past and shape our understanding of the world we live in"

Replace-Text " Unlike classical computers that rely on binary bits, quantum computers harness the enigmatic properties of quantum mechanics, enabling them to manipulate qubits" " It is a realm of gripping tales, both triumphant and tragic, revealing the echoes of human endeavors and the kaleidoscope of civilizations that have graced our planet"

Replace-Text " These qubits can exist in multiple states simultaneously, a phenomenon known as superposition, and become entangled, allowing for intricate correlations between them. This groundbreaking approach grants quantum computers an unparalleled computational prowess, opening up avenues for resolving currently intractable problems in diverse fields ranging from cryptography and optimization to simulations of complex systems" " History, as a discipline, endeavors to decipher the enigmas of our origins, traverse the intricate web of events that have shaped societies, and illuminate the lessons we can glean from the triumphs and missteps of those who came before us"

Replace-Text "As quantum computing advances, it promises to usher in a new era of technological progress" "In the tapestry of history, we find a mirror to ourselves, reflecting the struggles and aspirations that define the human experience"

Replace-Text " Its potential applications span across industries, spanning drug discovery, materials science, and financial modeling, to name a few" " As we delve into the chronicles of empires and civilizations, we uncover patterns of human behavior, cycles of progress and decline, and the threads that connect us to our ancestors"

Replace-Text " Quantum algorithms, such as Shor's algorithm for factoring large numbers and Grover's algorithm for searching unstructured databases, demonstrate the potential to revolutionize various computational tasks. By harnessing the quantum mechanical properties of particles, quantum computers can perform calculations that would take classical computers an impractical amount of time, thus unlocking new frontiers of scientific exploration and paving the way for groundbreaking discoveries" " History unveils the stories of remarkable individuals whose actions and decisions have shaped the course of events, reminding us of the profound impact that a single person can have"

Replace-Text "The development of quantum computers, however, is not without its challenges" "Furthermore, history is a vital tool for understanding the present, shedding light on the roots of current issues and helping us navigate the challenges of our time"

Replace-Text " Building and maintaining these intricate systems requires specialized expertise and infrastructure" " By studying the past, we can gain insights into the complexities of human nature, the dynamics of political and economic systems, and the forces that drive social change"

Replace-Text " Furthermore, the inherent fragility of quantum states poses a significant hurdle in preserving and manipulating information accurately. Despite these technical hurdles, significant progress has been made in recent years, with quantum computers achieving milestones such as demonstrating quantum supremacy and executing algorithms that cannot be efficiently executed on classical computers. The pursuit of quantum computing has sparked international scientific collaborations and investments from both academia and industry, propelling the field forward at an accelerated pace" " History provides a lens through which we can examine the motives of leaders, the aspirations of nations, and the interplay between different cultures"

# Summary paragraph
Replace-Text "Quantum computing represents a groundbreaking paradigm shift in computational capabilities, leveraging the principles of quantum mechanics to unlock unprecedented problem-solving potential" "History is a captivating subject that unveils the tapestry of human civilization and offers profound insights into the present"

Replace-Text " Its potential applications traverse diverse fields, ranging from cryptography and optimization to simulations of complex systems" " By exploring the annals of time, we uncover the intricacies of our origins, delve into the lessons of past triumphs and tribulations, and gain a deeper understanding of ourselves and the world we inhabit"

Replace-Text " While the field faces technical challenges, the rapid advancements in quantum computing hold immense promise for revolutionizing industries and expanding the boundaries of scientific inquiry. As we continue to delve into the realm of quantum computing, we stand at the cusp of a technological revolution poised to reimagine the very nature of computation and usher in a new era of scientific discovery" " History is a testament to the resilience, ingenuity, and creativity of humankind, inspiring us to reflect on our shared past and embrace the opportunities it presents to shape a better future"

# Add a new empty paragraph at the very end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
